$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3903.3962
$ws.Range("I64").Value = 3762.5
$ws.Range("J64").Value = 4020
$ws.Range("K64").Value = 3762.5
$ws.Range("L64").Value = 4020
$ws.Range("M64").Value = -3514.5
$ws.Range("N64").Value = -4516

$ws.Range("H67").Value = 3903.3962
$ws.Range("I67").Value = 3762.5
$ws.Range("J67").Value = 4020
$ws.Range("K67").Value = 3762.5
$ws.Range("L67").Value = 4020
$ws.Range("M67").Value = -2904.5
$ws.Range("N67").Value = -5736

$ws.Range("H125").Value = 1350.6666
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1350.6666
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 12155.9994
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -17075.9994

$ws.Range("H132").Value = 3595.3076
$ws.Range("I132").Value = 2809.9666
$ws.Range("J132").Value = 6213.1113
$ws.Range("K132").Value = 8429.899800000001
$ws.Range("L132").Value = 18639.3339
$ws.Range("M132").Value = -5899.899800000001
$ws.Range("N132").Value = -23699.3339

$ws.Range("H137").Value = 16669520
$ws.Range("I137").Value = 29413566
$ws.Range("J137").Value = 4231.231
$ws.Range("K137").Value = 88240698
$ws.Range("L137").Value = 12693.693
$ws.Range("M137").Value = -88238148
$ws.Range("N137").Value = -17793.693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11640599
$ws.Range("I32").Value = 16134896
$ws.Range("J32").Value = 30332.916
$ws.Range("K32").Value = 16134896
$ws.Range("L32").Value = 30332.916
$ws.Range("M32").Value = -16134609
$ws.Range("N32").Value = -30906.916

$ws.Range("H74").Value = 1979.081
$ws.Range("I74").Value = 1584.75
$ws.Range("J74").Value = 4502.8
$ws.Range("K74").Value = 1584.75
$ws.Range("L74").Value = 4502.8
$ws.Range("M74").Value = -710.75
$ws.Range("N74").Value = -6250.8

$ws.Range("H77").Value = 1979.081
$ws.Range("I77").Value = 1584.75
$ws.Range("J77").Value = 4502.8
$ws.Range("K77").Value = 7923.75
$ws.Range("L77").Value = 22514
$ws.Range("M77").Value = -3555.75
$ws.Range("N77").Value = -31250

$ws.Range("H117").Value = 43675
$ws.Range("J117").Value = 43675
$ws.Range("L117").Value = 43675
$ws.Range("N117").Value = -52853

$ws.Range("H122").Value = 2221.4
$ws.Range("I122").Value = 1400
$ws.Range("J122").Value = 2426.75
$ws.Range("K122").Value = 4200
$ws.Range("L122").Value = 7280.25
$ws.Range("M122").Value = -1750
$ws.Range("N122").Value = -12180.25

$ws.Range("H128").Value = 44500
$ws.Range("J128").Value = 44500
$ws.Range("L128").Value = 44500
$ws.Range("N128").Value = -54460

$ws.Range("H141").Value = 39800
$ws.Range("J141").Value = 39800
$ws.Range("L141").Value = 39800
$ws.Range("N141").Value = -50160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1239.5869
$ws.Range("I99").Value = 1040.3334
$ws.Range("K99").Value = 1040.3334
$ws.Range("M99").Value = 457.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1235.4615
$ws.Range("I122").Value = 951.875
$ws.Range("J122").Value = 1689.2
$ws.Range("K122").Value = 2855.625
$ws.Range("L122").Value = 5067.6
$ws.Range("M122").Value = -405.625
$ws.Range("N122").Value = -9967.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 555866.75
$ws.Range("I107").Value = 176.5
$ws.Range("J107").Value = 714635.4
$ws.Range("K107").Value = 529.5
$ws.Range("L107").Value = 2143906.2
$ws.Range("M107").Value = 1390.5
$ws.Range("N107").Value = -2147746.2

$ws.Range("H129").Value = 4260.05
$ws.Range("J129").Value = 5857.4
$ws.Range("L129").Value = 17572.2
$ws.Range("N129").Value = -27572.2

$ws.Range("H131").Value = 749.13043
$ws.Range("J131").Value = 924
$ws.Range("L131").Value = 2772
$ws.Range("N131").Value = -12852

$ws.Range("H132").Value = 4679618.5
$ws.Range("I132").Value = 1125
$ws.Range("J132").Value = 5230029.5
$ws.Range("K132").Value = 10125
$ws.Range("L132").Value = 47070265.5
$ws.Range("M132").Value = -7595
$ws.Range("N132").Value = -47075325.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 34000
$ws.Range("J95").Value = 34000
$ws.Range("L95").Value = 34000
$ws.Range("N95").Value = -39492

$ws.Range("H122").Value = 6252137
$ws.Range("I122").Value = 9092219
$ws.Range("J122").Value = 3955.4
$ws.Range("K122").Value = 27276657
$ws.Range("L122").Value = 11866.2
$ws.Range("M122").Value = -27274207
$ws.Range("N122").Value = -16766.2

$ws.Range("H126").Value = 14142878
$ws.Range("I126").Value = 11112677
$ws.Range("J126").Value = 27778784
$ws.Range("K126").Value = 33338031
$ws.Range("L126").Value = 83336352
$ws.Range("M126").Value = -33335561
$ws.Range("N126").Value = -83341292

$ws.Range("H138").Value = 26916
$ws.Range("J138").Value = 26916
$ws.Range("L138").Value = 26916
$ws.Range("N138").Value = -37196

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 921.8461
$ws.Range("I22").Value = 572
$ws.Range("K22").Value = 572
$ws.Range("M22").Value = -277

$ws.Range("H27").Value = 921.8461
$ws.Range("I27").Value = 572
$ws.Range("K27").Value = 572
$ws.Range("M27").Value = -465

$ws.Range("H46").Value = 877.2143
$ws.Range("I46").Value = 527.5
$ws.Range("J46").Value = 1139.5
$ws.Range("K46").Value = 527.5
$ws.Range("L46").Value = 1139.5
$ws.Range("M46").Value = -339.5
$ws.Range("N46").Value = -1515.5

$ws.Range("H87").Value = 32490
$ws.Range("J87").Value = 32490
$ws.Range("L87").Value = 32490
$ws.Range("N87").Value = -34736

$ws.Range("H90").Value = 32490
$ws.Range("J90").Value = 32490
$ws.Range("L90").Value = 97470
$ws.Range("N90").Value = -108702

$ws.Range("H122").Value = 2818.3333
$ws.Range("I122").Value = 2818.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8454.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6004.999899999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4133.3335
$ws.Range("I81").Value = 1500
$ws.Range("J81").Value = 4660
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 9320
$ws.Range("M81").Value = -1939
$ws.Range("N81").Value = -11442

$ws.Range("H84").Value = 4133.3335
$ws.Range("I84").Value = 1500
$ws.Range("J84").Value = 4660
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 46600
$ws.Range("M84").Value = -9696
$ws.Range("N84").Value = -57208

$ws.Range("H122").Value = 1960
$ws.Range("I122").Value = 1450.8823
$ws.Range("J122").Value = 2921.6667
$ws.Range("K122").Value = 4352.6469
$ws.Range("L122").Value = 8765.000100000001
$ws.Range("M122").Value = -1902.6469
$ws.Range("N122").Value = -13665.0001
